# The deck shipped with two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Integral" / "Red Violet" colour scheme (used by the slide master)
#   ppt/theme/theme2.xml -> "Office Theme" / "Office" colour scheme (used by the notes master)
#
# The authored change swaps the content of those two theme parts, i.e. the
# presentation's live design becomes the stock "Office" colour scheme.
# PowerPoint's object model exposes the live (slide-master-backed) theme's
# palette through SlideMaster.Theme.ThemeColorScheme, so re-point each of its
# twelve theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in that
# fixed index order) at the target "Office" theme's RGB values. MsoRGBType
# colours are packed as 0x00BBGGRR, so build each value from its hex triplet.

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0x000000   # dk1      #000000
$colorScheme.Item(2).RGB  = 0xFFFFFF   # lt1      #FFFFFF
$colorScheme.Item(3).RGB  = 0x6A5444   # dk2      #44546A
$colorScheme.Item(4).RGB  = 0xE6E6E7   # lt2      #E7E6E6
$colorScheme.Item(5).RGB  = 0xD59B5B   # accent1  #5B9BD5
$colorScheme.Item(6).RGB  = 0x317DED   # accent2  #ED7D31
$colorScheme.Item(7).RGB  = 0xA5A5A5   # accent3  #A5A5A5
$colorScheme.Item(8).RGB  = 0x00C0FF   # accent4  #FFC000
$colorScheme.Item(9).RGB  = 0xC47244   # accent5  #4472C4
$colorScheme.Item(10).RGB = 0x47AD70   # accent6  #70AD47
$colorScheme.Item(11).RGB = 0xC16305   # hlink    #0563C1
$colorScheme.Item(12).RGB = 0x724F95   # folHlink #954F72
